# ---------------------------------------------------------------------------
# Edit summary
# ---------------------------------------------------------------------------
# 1) Slide 5 contains a table (graphicFrame) whose table style is changed
#    from {A6C3EBDB-7FF8-4D8A-B248-44208946056D} to
#    {D7F16611-0BDE-4E95-B6DD-F5D61DBF93AE}.
#
# 2) The presentation's design/theme is switched from the "Integral" theme
#    (Red Violet color scheme) to the built-in "Office Theme" (Office color
#    scheme). In the COM object model the currently-applied theme's colors
#    are reached through Slide.ThemeColorScheme, whose 12 entries map to the
#    theme's dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink slots, in that order.
#    (PowerPoint's RGB color integers are packed as 0x00BBGGRR, same as the
#    VBA RGB() function.)
# ---------------------------------------------------------------------------

function ThemeRGB($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 5 --------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tbl = $tableShape.Table
$tbl.ApplyStyle("{D7F16611-0BDE-4E95-B6DD-F5D61DBF93AE}")

# --- 2) Apply the "Office Theme" colors -----------------------------------
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Item(1).RGB  = ThemeRGB 0x00 0x00 0x00   # dk1      000000
$tcs.Item(2).RGB  = ThemeRGB 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = ThemeRGB 0x44 0x54 0x6A   # dk2      44546A
$tcs.Item(4).RGB  = ThemeRGB 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Item(5).RGB  = ThemeRGB 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Item(6).RGB  = ThemeRGB 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Item(7).RGB  = ThemeRGB 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Item(8).RGB  = ThemeRGB 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Item(9).RGB  = ThemeRGB 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Item(10).RGB = ThemeRGB 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Item(11).RGB = ThemeRGB 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Item(12).RGB = ThemeRGB 0x95 0x4F 0x72   # folHlink 954F72
